$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.23512655154714679
$ws.Range("B1").Value = 0.23459094318251772
$ws.Range("A2").Value = -0.17190942954084854
$ws.Range("B2").Value = 0.17035052017249797
$ws.Range("A3").Value = -0.1206323867633472
$ws.Range("B3").Value = 0.12016268991811252
$ws.Range("A4").Value = -0.11216268995799261
$ws.Range("B4").Value = 0.11174233135237976
$ws.Range("A5").Value = -0.10874233137556999
$ws.Range("B5").Value = 0.10731770752222669
$ws.Range("A6").Value = -0.0080349826492849274
$ws.Range("B6").Value = 0.0079812981109981251
$ws.Range("A7").Value = 0.0020187018326498496
$ws.Range("B7").Value = -0.0020221677656904191
$ws.Range("A8").Value = 0.012022167709428988
$ws.Range("B8").Value = -0.012030433735799928
$ws.Range("A9").Value = 0.014030433711103463
$ws.Range("B9").Value = -0.014046010577732915
$ws.Range("A10").Value = 0.01604601055401389
$ws.Range("B10").Value = -0.01604760096786606
$ws.Range("A11").Value = 0.019047600940123033
$ws.Range("B11").Value = -0.019053926953730205
$ws.Range("A12").Value = 0.02255392692441216
$ws.Range("B12").Value = -0.022651396633118104
$ws.Range("A13").Value = 0.02615139660604715
$ws.Range("B13").Value = -0.026236675037913493
$ws.Range("A14").Value = 0.034236674994497562
$ws.Range("B14").Value = -0.034349200198218632
$ws.Range("A15").Value = 0.035349200182991147
$ws.Range("B15").Value = -0.035489311284076841
$ws.Range("A16").Value = 0.007375010413576355
$ws.Range("B16").Value = -0.0073871078998428708
$ws.Range("A17").Value = 0.0093871078830005672
$ws.Range("B17").Value = -0.0093966756540320162
$ws.Range("A18").Value = -0.056681513665864713
$ws.Range("B18").Value = 0.056582937530670563
$ws.Range("A19").Value = -0.052582937547676956
$ws.Range("B19").Value = 0.051870201888684075
$ws.Range("A20").Value = -0.04787020191098712
$ws.Range("B20").Value = 0.047670041686476594
$ws.Range("A21").Value = -0.043670041710003993
$ws.Range("B21").Value = 0.043365022910090367
$ws.Range("A22").Value = -0.045718133430986185
$ws.Range("B22").Value = 0.045502887884551591
$ws.Range("A23").Value = -0.040502887911449292
$ws.Range("B23").Value = 0.040099832359336673
$ws.Range("A24").Value = -0.020099832446226706
$ws.Range("B24").Value = 0.019999999911950894
$ws.Range("A25").Value = -0.097282724898549944
$ws.Range("B25").Value = 0.097155851931509218
$ws.Range("A26").Value = -0.094655851960002479
$ws.Range("B26").Value = 0.094492773624430981
$ws.Range("A27").Value = -0.091992773654737459
$ws.Range("B27").Value = 0.091028574610774182
$ws.Range("A28").Value = -0.089028574646655478
$ws.Range("B28").Value = 0.08836836742130405
$ws.Range("A29").Value = -0.081368367481835513
$ws.Range("B29").Value = 0.081176812309859514
$ws.Range("A30").Value = -0.021176812576976456
$ws.Range("B30").Value = 0.021024544740468887
$ws.Range("A31").Value = -0.014024544805554484
$ws.Range("B31").Value = 0.014001390855000295
$ws.Range("A32").Value = -0.004001390931898996
$ws.Range("B32").Value = 0.0039999999460693658
# Widen column B to match column A's width (target stored width 15.42578125).
# ColumnWidth is rounded by Excel to whole-pixel increments, so 14.67 is the
# input that lands on the closest achievable stored width (15.5).
$ws.Columns.Item(2).ColumnWidth = 14.67
